$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Interface")
$ws2 = $wb.Worksheets.Item(2)

# Insert a new row 2 ("All" module row) above the existing module rows,
# copying formatting from the row that is now row 3 (former row 2).
$ws1.Rows("2:2").Insert()
$ws1.Range("A3:E3").Copy()
$ws1.Range("A2:E2").PasteSpecial(-4122)

$ws1.Range("A2").Value2 = "FILI"
$ws1.Range("B2").Value2 = "Variable Annuity"
$ws1.Range("C2").Value2 = "All"
$ws1.Range("D2").Value2 = "FILI.xlsx"
$ws1.Range("E2").Value2 = "No"

# The previously-existing module rows (now rows 3-10) should read "No"
# for the Execute column; only the last row (Agent Module, row 11) keeps "Yes".
$ws1.Range("E3:E10").Value2 = "No"

# Extend the conditional formatting range to cover the new row.
$fc = $ws1.Range("D2:E10").FormatConditions
for ($i = 1; $i -le $fc.Count; $i++) {
    $fc.Item($i).ModifyAppliesToRange($ws1.Range("D2:E11"))
}

# Update the stored UI selection state to match the edited workbook.
$ws2.Activate()
$ws2.Range("A12:XFD14").Select()
$ws1.Activate()
$ws1.Range("E12").Select()
